# v0.16: HexagonalSquares are drawn on maze image.
#
# - "squares" sheet: the "F" marker becomes "T" and the "G" marker becomes "B";
#   cell C3 (previously blank "_") now also gets the "T" marker.
# - "grounds" sheet: cell D3 (previously blank "_") now gets the "ice" marker.
# - The active sheet/tab switches from "squares" to "grounds", with the
#   "grounds" sheet's selection moving from D2 to D4.

$wb = $excel.ActiveWorkbook

$squares = $wb.Worksheets.Item("squares")
$grounds = $wb.Worksheets.Item("grounds")

# --- Cell content changes --------------------------------------------------
$squares.Range("C2").Value = "T"
$squares.Range("C3").Value = "T"
$squares.Range("D3").Value = "B"

$grounds.Range("D3").Value = "ice"

# --- Active sheet / selection changes --------------------------------------
# "grounds" becomes the active sheet (was "squares"), selection moves to D4.
$grounds.Select()
$grounds.Range("D4").Select()
